# Ajustes de inventario: actualizar cantidades/costos calculados y
# registrar los nuevos productos agregados (filas 35-41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Actualiza cantidad (C), precio (D, solo donde cambia) y costo (F) ---
$updates = @(
    @{ Row = 2;  C = 24;  F = 53900 },
    @{ Row = 3;  C = 27;  F = 25600 },
    @{ Row = 4;  C = 9;   F = 48800 },
    @{ Row = 5;  C = 21;  F = 50300 },
    @{ Row = 6;  C = 40;  F = 23600 },
    @{ Row = 7;  C = 27;  D = 88000;  F = 42500 },
    @{ Row = 8;  C = 13;  F = 13800 },
    @{ Row = 9;  C = 6;   F = 45900 },
    @{ Row = 10; C = 34;  F = 21500 },
    @{ Row = 11; C = 30;  F = 38300 },
    @{ Row = 12; C = 59;  F = 13300 },
    @{ Row = 13; C = 24;  F = 53500 },
    @{ Row = 14; C = 9;   F = 106400 },
    @{ Row = 15; C = 49;  F = 33200 },
    @{ Row = 16; C = 8;   F = 63600 },
    @{ Row = 17; C = 13;  F = 65600 },
    @{ Row = 18; C = 14;  F = 28600 },
    @{ Row = 19; C = 35;  F = 55500 },
    @{ Row = 20; C = 17;  F = 17500 },
    @{ Row = 21; C = 0;   D = 125000; F = 50000 },
    @{ Row = 22; C = 2;   F = 28500 },
    @{ Row = 23; C = 153; D = 10000; F = 2163 },
    @{ Row = 24; C = 146; D = 10000; F = 2096 },
    @{ Row = 25; C = 16;  D = 12000; F = 3741 },
    @{ Row = 26; C = 34;  F = 2883 },
    @{ Row = 27; C = 1;   F = 60000 },
    @{ Row = 28; C = 3;   F = 30000 },
    @{ Row = 29; C = 38;  F = 3083 },
    @{ Row = 30; C = 6;   F = 6950 },
    @{ Row = 31; C = 24;  F = 1775 },
    @{ Row = 32; C = 46;  F = 1333 },
    @{ Row = 33; C = 48;  F = 625 },
    @{ Row = 34; C = 42;  F = 3750 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("C$r").Value = $u.C
    if ($u.ContainsKey('D')) {
        $ws.Range("D$r").Value = $u.D
    }
    $ws.Range("F$r").Value = $u.F
}

# --- Agrega los nuevos productos (filas 35-41) ---
$newRows = @(
    @{ Row = 35; A = "Cerveza Poker Lata";          B = "Cervezas";  C = 42; D = 12000; E = "2/13/2026"; F = 2941 },
    @{ Row = 36; A = "Cerveza Club Colombia Lata";  B = "Cervezas";  C = 67; D = 12000; E = "2/13/2026"; F = 3059 },
    @{ Row = 37; A = "Postobon y Pepsi";            B = "Gaseosas";  C = 47; D = 5000;  E = "2/13/2026"; F = 2083 },
    @{ Row = 38; A = "Coca Cola";                   B = "Gaseosas";  C = 10; D = 5000;  E = "2/13/2026"; F = 2500 },
    @{ Row = 39; A = "Cerveza Coronita";             B = "Cervezas";  C = 11; D = 10000; E = "2/13/2026"; F = 2833 },
    @{ Row = 40; A = "Vinos";                        B = "Otros";     C = 16; D = 60000; E = "2/13/2026"; F = 30000 },
    @{ Row = 41; A = "Bonfiest Bomba";               B = "Otros";     C = 27; D = 3000;  E = "2/13/2026"; F = 1100 }
)

foreach ($n in $newRows) {
    $r = $n.Row
    $ws.Range("A$r").Value = $n.A
    $ws.Range("B$r").Value = $n.B
    $ws.Range("C$r").Value = $n.C
    $ws.Range("D$r").Value = $n.D
    # Mantener la fecha como texto literal (igual que el resto de la columna),
    # no como fecha convertida a numero de serie.
    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $n.E
    $ws.Range("F$r").Value = $n.F
}
